# Delete the last slide (slide 15, the "Activity" slide) from the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
